$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue 2 4 '67.410.66'
Set-TextValue 2 5 '  +0.21%  '
Set-TextValue 3 4 '3.506.65'
Set-TextValue 3 5 '  -0.55%  '
Set-TextValue 4 5 '  -0.01%  '
Set-TextValue 5 4 '600.37'
Set-TextValue 5 5 '  +0.56%  '
Set-TextValue 6 4 '175.88'
Set-TextValue 6 5 '  +2.33%  '
Set-TextValue 7 5 '  +0.03%  '
Set-TextValue 8 5 '  -0.89%  '
Set-TextValue 9 5 '  -2.14%  '
Set-TextValue 10 4 '7.15'
Set-TextValue 10 5 '  -2.61%  '
Set-TextValue 11 5 '  -1.33%  '
Set-TextValue 12 4 '4.113.37'
Set-TextValue 12 5 '  -0.45%  '
Set-TextValue 13 4 '31.27'
Set-TextValue 13 5 '  +9.34%  '
Set-TextValue 15 4 '67.378.64'
Set-TextValue 16 5 '  -1.42%  '
Set-TextValue 17 4 '3.507.03'
Set-TextValue 17 5 '  -0.17%  '
Set-TextValue 18 5 '  -1.00%  '
Set-TextValue 19 4 '14.65'
Set-TextValue 19 5 '  +3.24%  '
Set-TextValue 20 4 '393.28'
Set-TextValue 20 5 '  -0.93%  '
Set-TextValue 21 4 '8.00'
Set-TextValue 21 5 '  -0.37%  '
Set-TextValue 22 4 '73.50'
Set-TextValue 22 5 '  -0.09%  '
Set-TextValue 23 4 '0.541'
Set-TextValue 23 5 '  +0.24%  '
Set-TextValue 24 5 '  -0.27%  '
Set-TextValue 25 5 '  -0.06%  '
Set-TextValue 26 4 '0.0000123'
Set-TextValue 26 5 '  -0.79%  '
Set-TextValue 27 4 '10.31'
Set-TextValue 27 5 '  +0.43%  '
Set-TextValue 28 5 '  -1.09%  '
Set-TextValue 29 5 '  -0.37%  '
Set-TextValue 30 5 '  -2.97%  '
Set-TextValue 31 5 '  -3.48%  '
Set-TextValue 32 5 '  -0.50%  '
Set-TextValue 33 4 '23.70'
Set-TextValue 33 5 '  -1.93%  '
Set-TextValue 34 4 '7.39'
Set-TextValue 34 5 '  -0.30%  '
Set-TextValue 35 5 '  +1.44%  '
Set-TextValue 36 4 '163.63'
Set-TextValue 36 5 '  -0.19%  '
Set-TextValue 37 5 '  +1.59%  '
Set-TextValue 38 5 '  -2.19%  '
Set-TextValue 39 4 '7.00'
Set-TextValue 39 5 '  +1.53%  '
Set-TextValue 40 4 '4.67'
Set-TextValue 40 5 '  -1.99%  '
Set-TextValue 41 5 '  -0.54%  '
Set-TextValue 42 4 '27.23'
Set-TextValue 42 5 '  +0.91%  '
Set-TextValue 43 5 '  -2.47%  '
Set-TextValue 44 4 '2.810.54'
Set-TextValue 44 5 '  +0.03%  '
Set-TextValue 45 4 '42.57'
Set-TextValue 45 5 '  -0.81%  '
Set-TextValue 46 5 '  -2.80%  '
Set-TextValue 47 5 '  -4.00%  '
Set-TextValue 48 4 '337.79'
Set-TextValue 48 5 '  -1.48%  '
Set-TextValue 49 5 '  -2.40%  '
Set-TextValue 50 4 '33.72'
Set-TextValue 50 5 '  +0.36%  '
Set-TextValue 51 5 '  -0.96%  '
